$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cells M1 and N1
$ws.Cells.Item(1,13).Value = "Energy_brochure(kJ)"
$ws.Cells.Item(1,14).Value = "Predicted_Ablation_Volume"

# Insert 4 new rows (bottom-to-top so row numbers stay stable during insertion)
$ws.Rows("29").Insert()
$ws.Rows("26").Insert()
$ws.Rows("23").Insert()
$ws.Rows("20").Insert()

# Fill new row 20
$ws.Cells.Item(20,1).Value = "Amica (Probe)"
$ws.Cells.Item(20,2).Value = "Probe"
$ws.Cells.Item(20,3).Value = 4
$ws.Cells.Item(20,4).Value = "ellipsoid"
$ws.Cells.Item(20,5).Value = 20
$ws.Cells.Item(20,6).Value = 900
$ws.Cells.Item(20,7).Value = "15.5 13.5 13.5"
$ws.Cells.Item(20,8).Value = "-11.5 0 0"
$ws.Cells.Item(20,9).Value = "1 0 0 0 1 0 0 0 1"
$ws.Cells.Item(20,10).Value = 19
$ws.Cells.Item(20,11).Value = 16.5
$ws.Cells.Item(20,12).Value = 16.5
$ws.Cells.Item(20,13).Value = 18
$ws.Cells.Item(20,14).Value = 21.6676

# Fill new row 24
$ws.Cells.Item(24,1).Value = "Amica (Probe)"
$ws.Cells.Item(24,2).Value = "Probe"
$ws.Cells.Item(24,3).Value = 8
$ws.Cells.Item(24,4).Value = "ellipsoid"
$ws.Cells.Item(24,5).Value = 40
$ws.Cells.Item(24,6).Value = 900
$ws.Cells.Item(24,7).Value = "24.5 18.0 18.0"
$ws.Cells.Item(24,8).Value = "-20.5 0 0"
$ws.Cells.Item(24,9).Value = "1 0 0 0 1 0 0 0 1"
$ws.Cells.Item(24,10).Value = 25
$ws.Cells.Item(24,11).Value = 21
$ws.Cells.Item(24,12).Value = 21
$ws.Cells.Item(24,13).Value = 36
$ws.Cells.Item(24,14).Value = 46.1814

# Fill new row 28
$ws.Cells.Item(28,1).Value = "Amica (Probe)"
$ws.Cells.Item(28,2).Value = "Probe"
$ws.Cells.Item(28,3).Value = 12
$ws.Cells.Item(28,4).Value = "ellipsoid"
$ws.Cells.Item(28,5).Value = 60
$ws.Cells.Item(28,6).Value = 900
$ws.Cells.Item(28,7).Value = "27.0 20.0 20.0"
$ws.Cells.Item(28,8).Value = "-22 0 0"
$ws.Cells.Item(28,9).Value = "1 0 0 0 1 0 0 0 1"
$ws.Cells.Item(28,10).Value = 30.5
$ws.Cells.Item(28,11).Value = 24
$ws.Cells.Item(28,12).Value = 24
$ws.Cells.Item(28,13).Value = 54
$ws.Cells.Item(28,14).Value = 73.5887

# Fill new row 32
$ws.Cells.Item(32,1).Value = "Amica (Probe)"
$ws.Cells.Item(32,2).Value = "Probe"
$ws.Cells.Item(32,3).Value = 16
$ws.Cells.Item(32,4).Value = "ellipsoid"
$ws.Cells.Item(32,5).Value = 80
$ws.Cells.Item(32,6).Value = 900
$ws.Cells.Item(32,7).Value = "33.0 23.0 23.0"
$ws.Cells.Item(32,8).Value = "-24 0 0"
$ws.Cells.Item(32,9).Value = "1 0 0 0 1 0 0 0 1"
$ws.Cells.Item(32,10).Value = 36.5
$ws.Cells.Item(32,11).Value = 27.5
$ws.Cells.Item(32,12).Value = 27.5
$ws.Cells.Item(32,13).Value = 72
$ws.Cells.Item(32,14).Value = 115.6237

# Set selection to N1 to match target view state
$ws.Range("N1").Select()